$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2, shifting existing rows 2:192 down to 3:193.
$ws.Rows(2).Insert()

# Populate the newly inserted row 2 with the new transaction data.
$ws.Range("E2").Value = "Deposit"
$ws.Range("N2").Value = "Crypto"
$ws.Range("P2").Value = "BTC"
$ws.Range("T2").Value = 10544.6102

# Restore view state: scroll back to top-left and select the header-ish block.
$ws.Application.ActiveWindow.ScrollColumn = 15
$ws.Range("E2:N3").Select()
